$wb = $excel.ActiveWorkbook

# Sheet 1: "PTP Away"
$ws1 = $wb.Worksheets.Item("PTP Away")

# Minor floating point refresh of existing Standard Deviation row values
$ws1.Cells.Item(3, 2).Value = 0.005919320605711683
$ws1.Cells.Item(3, 4).Value = 0.00230812780035107

# New "Mean" row
$ws1.Cells.Item(5, 1).Value = "Mean"
$ws1.Cells.Item(5, 2).Value = 0.02163277678097587
$ws1.Cells.Item(5, 3).Value = 0.03037585875043268
$ws1.Cells.Item(5, 4).Value = 0.008962291747283109

# Sheet 2: "PTP Close"
$ws2 = $wb.Worksheets.Item("PTP Close")

# Minor floating point refresh of existing Standard Deviation / Maximum row values
$ws2.Cells.Item(3, 2).Value = 0.000866440502589821
$ws2.Cells.Item(3, 3).Value = 0.0007782440124301294
$ws2.Cells.Item(3, 4).Value = 0.002557486753002339
$ws2.Cells.Item(4, 3).Value = 0.002422317488014691
$ws2.Cells.Item(4, 4).Value = 0.01084712176318814

# New "Mean" row
$ws2.Cells.Item(5, 1).Value = "Mean"
$ws2.Cells.Item(5, 2).Value = 0.003158025703890132
$ws2.Cells.Item(5, 3).Value = 0.0009597966896873986
$ws2.Cells.Item(5, 4).Value = 0.009795810698152472
